$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in shared string used by J1 header ("Residential" -> "Residentia")
$ws.Range("J1").Value = "Residentia"

# New column K: household expenditure change data ("houseex")
$ws.Range("K1").Value = "houseex"

$ws.Range("K2").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("K7").Value = -4.9292343457456536
$ws.Range("K8").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("K10").Value = -4.8203515127295331
$ws.Range("K11").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("K14").Value = -7.5196914542577344
$ws.Range("K15").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("K27").Value = -6.181305037989234
$ws.Range("K28").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("K30").Value = -6.3433217044006769
$ws.Range("K31").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("K35").Value = 0

# Match column K's best-fit width as closely as this runtime's width model allows
$ws.Columns.Item(11).ColumnWidth = 14

# Update selection to match the new view state
$ws.Range("K39").Select() | Out-Null
